# Generate Report for Handoff
#
# The "b.md" source file has completed a new localization handoff cycle.
# Update the Overview sheet's status for b.md, and update the per-locale
# (zh-cn / de-de) detail sheets with the new "Ready for handoff" status,
# the new handoff file name, and the new handoff datetime for b.md's row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: b.md row (row 3) status columns (zh-cn=B, de-de=C)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------------
# zh-cn detail sheet: b.md row (row 3)
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-03-10 16:29:44"

# Point the "Latest Handoff File" hyperlink at the new handoff file name
# (the link still resolves under the same handoff folder/commit as
# before; only the file name visible in the cell changes).
$zhcn.Hyperlinks.Add($zhcn.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/16e7fe78365e095bb5f1695cb04da118b27037c4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# de-de detail sheet: b.md row (row 3)
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("D3").Value = "2016-03-10 16:29:48"

$dede.Hyperlinks.Add($dede.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ce63d7c984156362b418d46dc37e035bd42aa38d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf") | Out-Null
